# Corrected Calibration and Ingest Sheets for Coastal Gliders
# - FLORT CC_angular_resolution -> 1.076
# - FLORT CC_scattering_angle   -> 124

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")
$ws.Activate()

# Row 4: CC_scattering_angle
$ws.Range("F4").Value = 124

# Row 6: CC_angular_resolution
$ws.Range("F6").Value = 1.076

# Leave the selection where the author left it when they saved the file
$ws.Range("E30").Select()
